$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update validation-rule text for "descrizione problema" (DD_Seg) and "quantità" (DD_QtOrd)
$ws.Range("B33").Value = "almeno 10 caratteri massimo 250 caratteri"
$ws.Range("B36").Value = "maggiore di 0 e minore uguale di 50"

# Remove the "immagine" attribute row from the DD_Prd section (old row 46),
# shifting "prezzo" up to take its place.
$ws.Rows("46:46").Delete()

# Update the view: scroll down and select B35 (as recorded in the saved file)
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("B35").Select()

# Restore the page setup recorded for the sheet
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
